$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) from 45302 to 45303 for existing rows 2-27
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 3).Value = 45303
}

# Row 27 gains an explicit row height (ht="15" customHeight="1")
$ws.Rows.Item(27).RowHeight = 15

# New row 28: A 1236-2024 (also gets explicit row height)
$ws.Cells.Item(28, 1).Value = "A 1236-2024"
$ws.Cells.Item(28, 2).Value = 45302
$ws.Cells.Item(28, 3).Value = 45303
$ws.Cells.Item(28, 4).Value = "OKÄNT"
$ws.Cells.Item(28, 5).Value = "OKÄNT"
$ws.Cells.Item(28, 7).Value = 0.7
$ws.Cells.Item(28, 8).Value = 0
$ws.Cells.Item(28, 9).Value = 0
$ws.Cells.Item(28, 10).Value = 0
$ws.Cells.Item(28, 11).Value = 0
$ws.Cells.Item(28, 12).Value = 0
$ws.Cells.Item(28, 13).Value = 0
$ws.Cells.Item(28, 14).Value = 0
$ws.Cells.Item(28, 15).Value = 0
$ws.Cells.Item(28, 16).Value = 0
$ws.Cells.Item(28, 17).Value = 0
$ws.Range("R28").WrapText = $true
$ws.Rows.Item(28).RowHeight = 15

# New row 29: A 1243-2024 (no explicit row height, like most of the new-row pair)
$ws.Cells.Item(29, 1).Value = "A 1243-2024"
$ws.Cells.Item(29, 2).Value = 45302
$ws.Cells.Item(29, 3).Value = 45303
$ws.Cells.Item(29, 4).Value = "OKÄNT"
$ws.Cells.Item(29, 5).Value = "OKÄNT"
$ws.Cells.Item(29, 7).Value = 6
$ws.Cells.Item(29, 8).Value = 0
$ws.Cells.Item(29, 9).Value = 0
$ws.Cells.Item(29, 10).Value = 0
$ws.Cells.Item(29, 11).Value = 0
$ws.Cells.Item(29, 12).Value = 0
$ws.Cells.Item(29, 13).Value = 0
$ws.Cells.Item(29, 14).Value = 0
$ws.Cells.Item(29, 15).Value = 0
$ws.Cells.Item(29, 16).Value = 0
$ws.Cells.Item(29, 17).Value = 0
$ws.Range("R29").WrapText = $true

# Match the date-formatted style (YYYY-MM-DD, style index 1) used by column B/C elsewhere
$ws.Range("B28:C29").NumberFormat = $ws.Range("B27:C27").NumberFormat
